$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the "A" column style (border/bold/center-top) used by rows 2-7
# down through the newly added rows 8-11, by copying the formatting of
# the existing A7 cell.
$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New coefficient table data: index, label, AR (C), DEC (D)
$rows = @(
    @(0, "X³",  [double]"-6.012381625503703e-12", [double]"-3.480997794492862e-11"),
    @(1, "X²Y", [double]"1.058815889446084e-11",  [double]"6.676038178035399e-13"),
    @(2, "XY²", [double]"-5.605361859614631e-11", [double]"1.49343710714519e-10"),
    @(3, "Y³",  [double]"1.555203461970007e-11",  [double]"-3.942116928602267e-11"),
    @(4, "X²",  [double]"1.570812355623091e-09",  [double]"7.98651376093351e-08"),
    @(5, "XY",  [double]"4.411993144177312e-08",  [double]"-1.501778878432668e-07"),
    @(6, "Y²",  [double]"1.076669423321491e-08",  [double]"-5.284675560796471e-08"),
    @(7, "X",   [double]"-1.680650016988583e-05", [double]"0.0002222532045222658"),
    @(8, "Y",   [double]"-0.0003017716375390833", [double]"6.887305019316614e-05"),
    @(9, "ind", [double]"293.1341250401416",      [double]"27.85148860818273")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
